$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add row 56 (Lp=55), wrapped text row, taller (ht=30)
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 3).Value = "ChatGPT, przy okazji eliminacji liczb dziesiętnych przy wczytywaniu nieobecności, pokazuje inną metodę zapisu do bazy z pliku Excel. Sprawdzić to, przeanalizować i zastodować. "
$ws.Cells.Item(56, 3).WrapText = $true
$ws.Cells.Item(56, 4).Value = 0
$ws.Rows.Item(56).RowHeight = 30

# Add row 57 (Lp=56)
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "Nieobecności"
$ws.Cells.Item(57, 3).Value = "Filtrowanie danych czytanych z pliku i eliminacja liczb dziesiętnych. Zapis tylko całkowitych liczb"
$ws.Cells.Item(57, 3).WrapText = $true
$ws.Cells.Item(57, 4).Value = 0

# Selection and frozen pane view
$ws.Range("C57").Select()
